$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Progress value in E3 changes from 0.8 (80%) to 0 (0%)
$ws.Range("E3").Value = 0

# New "Subject" entry added in F3
$ws.Range("F3").Value = "Decision Tree"

# Rows 4-17 get a Progress value of 0%, using the same percentage
# number format already applied to E2:E3
$ws.Range("E4:E17").NumberFormat = $ws.Range("E3").NumberFormat
$ws.Range("E4:E17").Value = 0

# Active cell/selection moves from F2 to D3
$ws.Range("D3").Select()
